# Fruta / hortaliza, semanal
# Insert 6 new weekly price rows (Comercializadora del Agro de Limari - Nectarin)
# above the existing "August Red" block, pushing the existing rows 169-178
# down to 175-184, and populate the new rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at 169; existing rows 169-178 shift down to 175-184.
$ws.Rows("169:174").Insert()

# Columns that stay constant across this whole producer/product block.
$colA = 2
$colB = "Comercializadora del Agro de Limarí"
$colC = "Coquimbo"
$colE = 4
$colF = "Fruta"
$colG = 100103
$colH = "Frutos de hueso (carozo)"
$colI = 100103006
$colJ = "Nectarín"
$colR = "Región de O'Higgins"

# New rows data: Date(serial), Variedad, Calidad, Cantidad, Min, Max, Prom, Unidad, S, T
$newRows = @(
    @{ Row = 169; D = 44615; K = "August Red"; L = "Primera";  M = 20; N = 320000; O = 330000; P = 325000; Q = "`$/bins (420 kilos)"; S = 774; T = 420 },
    @{ Row = 170; D = 44615; K = "August Red"; L = "Segunda";  M = 20; N = 270000; O = 280000; P = 275000; Q = "`$/bins (420 kilos)"; S = 655; T = 420 },
    @{ Row = 171; D = 44615; K = "June Pearl"; L = "Especial"; M = 16; N = 360000; O = 370000; P = 365000; Q = "`$/bins (420 kilos)"; S = 869; T = 420 },
    @{ Row = 172; D = 44615; K = "June Pearl"; L = "Primera";  M = 20; N = 330000; O = 340000; P = 335000; Q = "`$/bins (420 kilos)"; S = 798; T = 420 },
    @{ Row = 173; D = 44615; K = "Venus";      L = "Especial"; M = 28; N = 320000; O = 330000; P = 323571; Q = "`$/bins (420 kilos)"; S = 770; T = 420 },
    @{ Row = 174; D = 44615; K = "Venus";      L = "Primera";  M = 20; N = 300000; O = 310000; P = 305000; Q = "`$/bins (420 kilos)"; S = 726; T = 420 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
    $ws.Cells.Item($row, 10).Value = $colJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $colR
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
